# "update employer quick registration"
#
# The quick-registration sample row (row 4) on "Shee1" held a work-history
# entry for "Meltshop ny" whose address was abbreviated ("111 Fulton St").
# Spell it out ("111 Fulton Street"), and drop the stray row of empty
# placeholder whitespace that used to sit just below the sample data in H5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- H4: "...111 Fulton St, ..." -> "...111 Fulton Street, ..." ----------
$ws.Range("H4").Value = "Meltshop ny;111 Fulton Street, New York, NY 10038, United States;Marketing Manager;March;2012;;;true"

# --- H5 used to hold a placeholder string of blank spaces; remove it -----
$ws.Range("H5").ClearContents()

# --- column widths were hand-tuned after the edit -------------------------
$ws.Columns.Item(8).ColumnWidth = 55        # H
$ws.Columns.Item(10).ColumnWidth = 63.666   # J
$ws.Columns.Item(11).ColumnWidth = 27.333   # K
$ws.Columns.Item(12).ColumnWidth = 21.333   # L

# --- selection / scroll position left where the edit was made ------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 9          # column I -> topLeftCell "I1"
$ws.Range("J4").Select() | Out-Null

# --- window size/position as last saved -----------------------------------
$win.Left = 2625
$win.Top = 2430
$win.Width = 14790
$win.Height = 6405
